# Menu expansion: append five new pizza rows (name / type / price) below
# the existing price table on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 6;  Name = "Mushroom";       Type = "VEGETERIAN"; Price = 245 },
    @{ Row = 7;  Name = "Cheese";         Type = "VEGETERIAN"; Price = 354 },
    @{ Row = 8;  Name = "calamari";       Type = "SEA";         Price = 454 },
    @{ Row = 9;  Name = "DoubleCheese";   Type = "VEGETERIAN"; Price = 453 },
    @{ Row = 10; Name = "DoubleMushroom"; Type = "VEGETERIAN"; Price = 200 }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.Name
    $ws.Cells.Item($r.Row, 3).Value = $r.Type
    $ws.Cells.Item($r.Row, 4).Value = $r.Price
}

[void]$ws.Range("A10").Select()

# Best-effort re-creation of the two differential-style ("dxf") records that
# ship with the custom "MySqlDefault" table style in the target workbook
# (bold/shaded header row + plain "whole table" body). Excel only persists
# dxf records once something references them, so we stamp them onto the
# table range via a pair of transient conditional-format rules and then
# remove the rules again -- same way Excel itself leaves "orphaned" dxfs
# behind in styles.xml after a format is defined and the rule that used it
# is removed.
$styleRange = $ws.Range("B2:D10")

$headerDxf = $styleRange.FormatConditions.Add(2, 5, "0")
$headerDxf.Font.Bold = $true
$headerDxf.Interior.Color = 14145495

$wholeTableDxf = $styleRange.FormatConditions.Add(2, 5, "0")
$wholeTableDxf.Font.Bold = $false
$wholeTableDxf.Interior.Pattern = -4142
$wholeTableDxf.Interior.ColorIndex = -4142

$styleRange.FormatConditions.Delete()

Write-Output "done"

